# Add team record (Wins/Losses/Ties) columns to the BAL_1991 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row currently used by the sheet (rows 2..44 hold player data).
$lastRow = 44

# --- Header row (row 1): new column headers in AD1:AF1 -----------------
# First clone the existing header style (bold, centered, bordered) from the
# neighbouring "Unnamed: 28" header cell (AC1) onto the three new header
# cells, then overwrite their text so the format matches the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2..44): team record repeated for every player ----------
$wins = 67
$losses = 95
$ties = 0

$ws.Range("AD2:AD$lastRow").Value = $wins
$ws.Range("AE2:AE$lastRow").Value = $losses
$ws.Range("AF2:AF$lastRow").Value = $ties

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
